$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -4685
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3908
$ws.Range("N79").ClearContents()

$ws.Range("H86").Value = 7846.1665
$ws.Range("I86").Value = 6766.6665
$ws.Range("J86").Value = 8925.666999999999
$ws.Range("K86").Value = 6766.6665
$ws.Range("L86").Value = 8925.666999999999
$ws.Range("M86").Value = -5643.6665
$ws.Range("N86").Value = -11171.667

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H89").Value = 7846.1665
$ws.Range("I89").Value = 6766.6665
$ws.Range("J89").Value = 8925.666999999999
$ws.Range("K89").Value = 33833.3325
$ws.Range("L89").Value = 44628.335
$ws.Range("M89").Value = -28217.3325
$ws.Range("N89").Value = -55860.335

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H135").Value = 536.8182
$ws.Range("I135").Value = 280.5
$ws.Range("J135").Value = 1220.3334
$ws.Range("K135").Value = 2524.5
$ws.Range("L135").Value = 10983.0006
$ws.Range("M135").Value = 10.5
$ws.Range("N135").Value = -16053.0006

$ws.Range("H137").Value = 1747
$ws.Range("I137").Value = 1725.6
$ws.Range("J137").Value = 1782.6666
$ws.Range("K137").Value = 5176.799999999999
$ws.Range("L137").Value = 5347.9998
$ws.Range("M137").Value = -2626.799999999999
$ws.Range("N137").Value = -10447.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3421.0344
$ws.Range("J32").Value = 8013.4
$ws.Range("L32").Value = 8013.4
$ws.Range("N32").Value = -8587.4

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 1037.7059
$ws.Range("I74").Value = 1037.7059
$ws.Range("K74").Value = 1037.7059
$ws.Range("M74").Value = -163.7058999999999

$ws.Range("H77").Value = 1037.7059
$ws.Range("I77").Value = 1037.7059
$ws.Range("K77").Value = 5188.5295
$ws.Range("M77").Value = -820.5294999999996

$ws.Range("H102").Value = 2271.389
$ws.Range("I102").Value = 2257.9412
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2257.9412
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -635.9412000000002
$ws.Range("N102").Value = -5744

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 766.2727
$ws.Range("I64").Value = 887
$ws.Range("J64").Value = 739.44446
$ws.Range("K64").Value = 887
$ws.Range("L64").Value = 739.44446
$ws.Range("M64").Value = -662
$ws.Range("N64").Value = -1189.44446

$ws.Range("H67").Value = 766.2727
$ws.Range("I67").Value = 887
$ws.Range("J67").Value = 739.44446
$ws.Range("K67").Value = 887
$ws.Range("L67").Value = 739.44446
$ws.Range("M67").Value = -107
$ws.Range("N67").Value = -2299.44446

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1860.4286
$ws.Range("I31").Value = 1904.6
$ws.Range("K31").Value = 1904.6
$ws.Range("M31").Value = -1609.6

$ws.Range("H34").Value = 1860.4286
$ws.Range("I34").Value = 1904.6
$ws.Range("K34").Value = 1904.6
$ws.Range("M34").Value = -1702.6

$ws.Range("H132").Value = 2647.261
$ws.Range("I132").Value = 2741.5715
$ws.Range("J132").Value = 1657
$ws.Range("K132").Value = 8224.7145
$ws.Range("L132").Value = 4971
$ws.Range("M132").Value = -5694.7145
$ws.Range("N132").Value = -10031

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 123.5
$ws.Range("I40").Value = 96.5
$ws.Range("K40").Value = 386
$ws.Range("M40").Value = -317

$ws.Range("H74").Value = 8400
$ws.Range("J74").Value = 8400
$ws.Range("L74").Value = 25200
$ws.Range("N74").Value = -27322

$ws.Range("H77").Value = 8400
$ws.Range("J77").Value = 8400
$ws.Range("L77").Value = 75600
$ws.Range("N77").Value = -86208

$ws.Range("H131").Value = 945.5454999999999
$ws.Range("I131").Value = 867.75
$ws.Range("J131").Value = 990
$ws.Range("K131").Value = 2603.25
$ws.Range("L131").Value = 2970
$ws.Range("M131").Value = 2436.75
$ws.Range("N131").Value = -13050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 963.6923
$ws.Range("I97").Value = 669
$ws.Range("K97").Value = 669
$ws.Range("M97").Value = -173

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 666.5
$ws.Range("J22").Value = 699.8
$ws.Range("L22").Value = 699.8
$ws.Range("N22").Value = -1289.8

$ws.Range("H27").Value = 666.5
$ws.Range("J27").Value = 699.8
$ws.Range("L27").Value = 699.8
$ws.Range("N27").Value = -913.8

$ws.Range("H55").Value = 373.25
$ws.Range("I55").Value = 61.285713
$ws.Range("K55").Value = 61.285713
$ws.Range("M55").Value = 111.714287

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H136").Value = 4500
$ws.Range("I136").Value = 4500
$ws.Range("K136").Value = 13500
$ws.Range("M136").Value = -10950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1900
$ws.Range("I96").Value = 1900
$ws.Range("K96").Value = 1900
$ws.Range("M96").Value = -527

$ws.Range("H136").Value = 1636.5454
$ws.Range("I136").Value = 1649.7
$ws.Range("K136").Value = 4949.1
$ws.Range("M136").Value = -2399.1
